# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) on the per-job Leve Profit sheets, pulling fresh Universalis-style values.
$wb = $excel.ActiveWorkbook

# ALC!row98 - The Dotted Line / Enchanted Durium Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 992
$ws.Range("I98").Value = 681.25
$ws.Range("J98").Value = 1737.8
$ws.Range("K98").Value = 681.25
$ws.Range("L98").Value = 1737.8
$ws.Range("M98").Value = 816.75
$ws.Range("N98").Value = -4733.8

# ALC!row112 - Making Ends Meet / Superior Spiritbond Potion
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 6183.3335
$ws.Range("I112").Value = 900
$ws.Range("J112").Value = 6365.517
$ws.Range("K112").Value = 2700
$ws.Range("L112").Value = 19096.551
$ws.Range("M112").Value = -1592
$ws.Range("N112").Value = -21312.551

# ALC!row122 - Wishful Inking / Enchanted High Durium Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 992
$ws.Range("I122").Value = 681.25
$ws.Range("J122").Value = 1737.8
$ws.Range("K122").Value = 2043.75
$ws.Range("L122").Value = 5213.4
$ws.Range("M122").Value = 406.25
$ws.Range("N122").Value = -10113.4

# ALC!row134 - Binding Spells / Crocodileskin Index
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 70140
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 70140
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 70140
$ws.Range("N134").Value = -80280

# ALC!row135 - For Tired Minds / Grade 1 Gemsap of Intelligence
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 616.2105
$ws.Range("I135").Value = 567.1429000000001
$ws.Range("J135").Value = 1188.6666
$ws.Range("K135").Value = 5104.2861
$ws.Range("L135").Value = 10697.9994
$ws.Range("M135").Value = -2569.2861
$ws.Range("N135").Value = -15767.9994

# ALC!row138 - All-night Crafting / Cunning Craftsman's Tisane
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1595.21
$ws.Range("I138").Value = 645.3214
$ws.Range("J138").Value = 1964.6111
$ws.Range("K138").Value = 1935.9642
$ws.Range("L138").Value = 5893.8333
$ws.Range("M138").Value = 3204.0358
$ws.Range("N138").Value = -16173.8333

# ARM!row32 - Ingot We Trust / Steel Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 756815.9
$ws.Range("I32").Value = 849425.6
$ws.Range("J32").Value = 25198.7
$ws.Range("K32").Value = 849425.6
$ws.Range("L32").Value = 25198.7
$ws.Range("M32").Value = -849138.6
$ws.Range("N32").Value = -25772.7

# ARM!row122 - Haste for High Durium / High Durium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1647.2142
$ws.Range("I122").Value = 1355.1818
$ws.Range("J122").Value = 2718
$ws.Range("K122").Value = 4065.5454
$ws.Range("L122").Value = 8154
$ws.Range("M122").Value = -1615.5454
$ws.Range("N122").Value = -13054

# ARM!row134 - Brace for More Vambraces / Ruthenium Vambraces of Maiming
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 40429
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 40429
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 40429
$ws.Range("N134").Value = -50569

# CRP!row31 - Wall Not Found / Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4094.7207
$ws.Range("I31").Value = 1264.6052
$ws.Range("J31").Value = 7679.533
$ws.Range("K31").Value = 1264.6052
$ws.Range("L31").Value = 7679.533
$ws.Range("M31").Value = -969.6052
$ws.Range("N31").Value = -8269.532999999999

# CRP!row34 - Armoires of the Rich and Famous / Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4094.7207
$ws.Range("I34").Value = 1264.6052
$ws.Range("J34").Value = 7679.533
$ws.Range("K34").Value = 1264.6052
$ws.Range("L34").Value = 7679.533
$ws.Range("M34").Value = -1062.6052
$ws.Range("N34").Value = -8083.533

# CRP!row94 - Beech, Please / Beech Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1306.9412
$ws.Range("I94").Value = 850
$ws.Range("J94").Value = 1367.8667
$ws.Range("K94").Value = 850
$ws.Range("L94").Value = 1367.8667
$ws.Range("M94").Value = -399
$ws.Range("N94").Value = -2269.8667

# CRP!row99 - O Pine / Pine Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1793.7727
$ws.Range("I99").Value = 1366.5
$ws.Range("J99").Value = 1888.7222
$ws.Range("K99").Value = 1366.5
$ws.Range("L99").Value = 1888.7222
$ws.Range("M99").Value = 131.5
$ws.Range("N99").Value = -4884.7222

# CRP!row116 - The Right Tool for the Job / Sandteak Rod
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 80000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 80000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 80000
$ws.Range("N116").Value = -89178

# CRP!row122 - Timber of Tenkonto / Horse Chestnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1702.2759
$ws.Range("I122").Value = 1497.2222
$ws.Range("J122").Value = 1794.55
$ws.Range("K122").Value = 4491.6666
$ws.Range("L122").Value = 5383.65
$ws.Range("M122").Value = -2041.6666
$ws.Range("N122").Value = -10283.65

# CRP!row126 - A Better Conductor / Red Pine Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1793.7727
$ws.Range("I126").Value = 1366.5
$ws.Range("J126").Value = 1888.7222
$ws.Range("K126").Value = 4099.5
$ws.Range("L126").Value = 5666.1666
$ws.Range("M126").Value = -1629.5
$ws.Range("N126").Value = -10606.1666

# CUL!row17 - Chew the Fat / Grilled Dodo
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 487.8889
$ws.Range("I17").Value = 365.16666
$ws.Range("J17").Value = 733.3333
$ws.Range("K17").Value = 1095.49998
$ws.Range("L17").Value = 2199.9999
$ws.Range("M17").Value = -926.4999800000001
$ws.Range("N17").Value = -2537.9999

# CUL!row113 - Can't Eat Just One / Night Vinegar
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 507.86667
$ws.Range("I113").Value = 506.13635
$ws.Range("J113").Value = 509.52173
$ws.Range("K113").Value = 1518.40905
$ws.Range("L113").Value = 1528.56519
$ws.Range("M113").Value = 651.59095
$ws.Range("N113").Value = -5868.56519

# GSM!row102 - Put the Metal to the Peddle / Durium Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1699.3405
$ws.Range("I102").Value = 1593
$ws.Range("J102").Value = 2009.5
$ws.Range("K102").Value = 1593
$ws.Range("L102").Value = 2009.5
$ws.Range("M102").Value = 29
$ws.Range("N102").Value = -5253.5

# GSM!row118 - A Magnanimous Refrain / Triplite Earrings of Casting
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 27560
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 27560
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 27560
$ws.Range("N118").Value = -30874

# GSM!row132 - On Board for Lar / Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2255.9058
$ws.Range("I132").Value = 1630.683
$ws.Range("J132").Value = 4392.0835
$ws.Range("K132").Value = 4892.049
$ws.Range("L132").Value = 13176.2505
$ws.Range("M132").Value = -2362.049
$ws.Range("N132").Value = -18236.2505

# LTW!row122 - Hell on Leather / Gaja Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3528.1365
$ws.Range("I122").Value = 3003.0908
$ws.Range("J122").Value = 4053.182
$ws.Range("K122").Value = 9009.2724
$ws.Range("L122").Value = 12159.546
$ws.Range("M122").Value = -6559.2724
$ws.Range("N122").Value = -17059.546

# WVR!row46 - Crunching the Numbers / Linen Hat
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 106107
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 106107
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 106107
$ws.Range("N46").Value = -106569

# WVR!row117 - The Hunt Continues / Ovim Wool Muffed Met of Casting
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 23574.25
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 23574.25
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 23574.25
$ws.Range("N117").Value = -32752.25

# WVR!row132 - Comfy Cabins / Snow Cotton Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6342646
$ws.Range("I132").Value = 2266.724
$ws.Range("J132").Value = 17158586
$ws.Range("K132").Value = 6800.172
$ws.Range("L132").Value = 51475758
$ws.Range("M132").Value = -4270.172
$ws.Range("N132").Value = -51480818

# WVR!row134 - Cloth for Canvas / Mountain Linen
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 106107
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 106107
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 318321
$ws.Range("N134").Value = -323391

# WVR!row136 - Weaving the Envelope / Sarcenet Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2011.51
$ws.Range("I136").Value = 1972.079
$ws.Range("J136").Value = 2136.375
$ws.Range("K136").Value = 5916.237
$ws.Range("L136").Value = 6409.125
$ws.Range("M136").Value = -3366.237
$ws.Range("N136").Value = -11509.125
